# Update the "user" column (E) values so that each row's User@<hex> identity
# hash reflects the newly created admin user/profile objects, as described by
# the commit "ADD: permisions admin user and profile".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "lv.venta.models.users.User@608b3a7"
$ws.Range("E3").Value = "lv.venta.models.users.User@577919f9"
$ws.Range("E4").Value = "lv.venta.models.users.User@437a9c25"
$ws.Range("E5").Value = "lv.venta.models.users.User@71b5abca"
$ws.Range("E6").Value = "lv.venta.models.users.User@6d5c7846"
$ws.Range("E7").Value = "lv.venta.models.users.User@14a7d2b1"
$ws.Range("E8").Value = "lv.venta.models.users.User@7df7ce8f"
